# Add two new columns (I: "I0", J: "IF") to the stats sheet, mirroring the
# existing header/style conventions already used by column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold / centered / bordered) used by the
# other header cells (e.g. H1) by copying its format onto I1:J1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (rows 2-29) ---------------------------------------------
$rows = @(
    @(2, 5, 7),
    @(3, 5, 6),
    @(4, 1, 5),
    @(5, 1, 5),
    @(6, 1, 6),
    @(7, 1, 6),
    @(8, 1, 6),
    @(9, 1, 5),
    @(10, 1, 4),
    @(11, 1, 6),
    @(12, 1, 6),
    @(13, 1, 4),
    @(14, 1, 5),
    @(15, 1, 2),
    @(16, 1, 6),
    @(17, 1, 6),
    @(18, 1, 6),
    @(19, 1, 5),
    @(20, 1, 6),
    @(21, 1, 6),
    @(22, 1, 5),
    @(23, 1, 5),
    @(24, 1, 5),
    @(25, 1, 5),
    @(26, 1, 4),
    @(27, 1, 3),
    @(28, 4, 5),
    @(29, 1, 2)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $iVal = $r[1]
    $jVal = $r[2]
    $ws.Cells.Item($rowNum, 9).Value = $iVal   # column I
    $ws.Cells.Item($rowNum, 10).Value = $jVal  # column J
}
